$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting existing rows 35..196 down to 36..197
$ws.Rows("35").Insert()

# Populate the newly inserted row 35 with the new record's data
$ws.Range("A35").Value = 8
$ws.Range("B35").Value = "Terminal La Palmera de La Serena"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 44547
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 100112012
$ws.Range("G35").Value = "Espinaca"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 3320
$ws.Range("K35").Value = 400
$ws.Range("L35").Value = 500
$ws.Range("M35").Value = 450
$ws.Range("N35").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O35").Value = "Provincia del Elquí"
$ws.Range("P35").Value = 900
$ws.Range("Q35").Value = 0.5
$ws.Range("R35").Value = "Hortaliza"
